# Update "Australia ALeague" odds data sheet.
# 1) Rows 104 and 105 had their match data (columns B:AC) swapped back
#    to the correct fixture order (id column A stays put).
# 2) Rows 112 and 113 likewise have columns B:AC swapped.
# 3) Rows 115-119 get individual odds cells refreshed with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap helper cells B (col 2) through AC (col 29) between two rows.
# NOTE: named parameter binding for script functions is unreliable in this
# runtime, so pass arguments positionally and avoid named parameters.
$rowPairs = @(
    @(104, 105),
    @(112, 113)
)

foreach ($pair in $rowPairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]
    for ($col = 2; $col -le 29; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# Row 115 - refreshed odds
$ws.Range("R115").Value = 1.86
$ws.Range("S115").Value = 2.04

# Row 116 - refreshed odds
$ws.Range("N116").Value = 3.4
$ws.Range("P116").Value = 2
$ws.Range("Q116").Value = 0.5
$ws.Range("R116").Value = 1.88
$ws.Range("S116").Value = 2.02
$ws.Range("U116").Value = 1.875
$ws.Range("V116").Value = 1.975

# Row 117 - refreshed odds
$ws.Range("R117").Value = 1.9
$ws.Range("S117").Value = 2
$ws.Range("U117").Value = 1.925
$ws.Range("V117").Value = 1.925

# Row 118 - refreshed odds
$ws.Range("N118").Value = 2.2
$ws.Range("O118").Value = 3.5
$ws.Range("P118").Value = 3.1
$ws.Range("Q118").Value = -0.25
$ws.Range("R118").Value = 1.97
$ws.Range("S118").Value = 1.93

# Row 119 - refreshed odds
$ws.Range("P119").Value = 2.2
$ws.Range("R119").Value = 1.92
$ws.Range("S119").Value = 1.98
$ws.Range("T119").Value = 3
$ws.Range("U119").Value = 2
$ws.Range("V119").Value = 1.85
